# Update the questionnaire title and the first three hardware-feature
# questions from the APS/PDW-specific wording to the generic MPP/Appliance
# wording, per the "Updated Deployment utility to allow mutiple variables"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Title
$ws.Range("B1").Value = "MPP Questionnaire"

# Hardware Features questions
$ws.Range("B3").Value = "Make and Model of your Appliance for each Environment? (CPUs, GPUs, memory)"
$ws.Range("B4").Value = "List the storage space, Memory, Cores and/or FGPAs for each environment?"
$ws.Range("B5").Value = "Is there a Loading Server?  If so, is it connected to the Appliance?"

$ws.Range("B6").Select()
